$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C23: 4.5 -> 4
$ws.Range("C23").Value2 = 4

# D23: add "Dash plotly" label, matching style of D21/D22 (center aligned)
$ws.Range("D23").Value2 = "Dash plotly"
$ws.Range("D23").HorizontalAlignment = -4108

# Row 24: Lunes / 2-dec / 1 / Dash plotly
$ws.Range("A24").Value2 = "Lunes"
$ws.Range("B24").Value2 = "2-dec"
$ws.Range("C24").Value2 = 1
$ws.Range("D24").Value2 = "Dash plotly"
$ws.Range("D24").HorizontalAlignment = -4108
$ws.Rows.Item(24).RowHeight = 13.8

# Row 25: martes / 3-dic / 2 / Dash plotly
$ws.Range("A25").Value2 = "martes"
$ws.Range("B25").Value2 = "3-dic"
$ws.Range("C25").Value2 = 2
$ws.Range("D25").Value2 = "Dash plotly"
$ws.Range("D25").HorizontalAlignment = -4108
$ws.Rows.Item(25).RowHeight = 13.8

# Update selection to D16
[void]$ws.Range("D16").Select()
